{"js": "// Change 1: \"A cool quote by Dijkstra:\" -> \"A cool quote by Archana:\"\n// Only the name changes; replace just \"Dijkstra\" so the surrounding\n// (unformatted) text is left untouched.\nconst nameResults = context.document.body.search(\"Dijkstra\", { matchCase: true });\nnameResults.load(\"text\");\nawait context.sync();\n\nif (nameResults.items.length > 0) {\n  nameResults.items[0].insertText(\"Archana\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Change 2: the quoted sentence run becomes \"What's stopping you?\" while the\n// closing curly quote (\") that ends the paragraph is preserved.\nconst quoteResults = context.document.body.search(\n  \"Computer science is no more about computers than astronomy is about telescopes.\",\n  { matchCase: true }\n);\nquoteResults.load(\"text\");\nawait context.sync();\n\nif (quoteResults.items.length > 0) {\n  quoteResults.items[0].insertText(\"What\\u2019s stopping you?\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Change 1: \"A cool quote by Dijkstra:\" -> \"A cool quote by Archana:\"\n# Replace only the name so the rest of the (unformatted) sentence is left as-is.\n$rngName = $d.Content\n$rngName.Find.ClearFormatting()\n$rngName.Find.Text = \"Dijkstra\"\n$rngName.Find.MatchCase = $true\n$rngName.Find.Execute() | Out-Null\nif ($rngName.Find.Found) {\n    $rngName.Text = \"Archana\"\n}\n\n# Change 2: the quoted sentence becomes \"What's stopping you?\" while the\n# closing curly quote (\") that ends the paragraph is kept untouched.\n$rngQuote = $d.Content\n$rngQuote.Find.ClearFormatting()\n$rngQuote.Find.Text = \"Computer science is no more about computers than astronomy is about telescopes.\"\n$rngQuote.Find.MatchCase = $true\n$rngQuote.Find.Execute() | Out-Null\nif ($rngQuote.Find.Found) {\n    $rngQuote.Text = \"What\" + [char]0x2019 + \"s stopping you?\"\n}\n"}
